$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032690091378365
$ws.Range("D2").Value = 1.037131897358414
$ws.Range("E2").Value = 1.032103089698571
$ws.Range("I2").Value = 1.037664730688466
$ws.Range("J2").Value = 1.03781892621172
$ws.Range("K2").Value = 1.039923656061474
$ws.Range("L2").Value = 1.034909287050578
$ws.Range("N2").Value = 1.016579996623338

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033630736533736
$ws.Range("D3").Value = 1.037855268837407
$ws.Range("E3").Value = 1.032901549108925
$ws.Range("I3").Value = 1.037926604017113
$ws.Range("J3").Value = 1.03840201669346
$ws.Range("K3").Value = 1.04045700096332
$ws.Range("L3").Value = 1.035516476323678
$ws.Range("N3").Value = 1.016773882354107

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034239464704153
$ws.Range("D4").Value = 1.038323270142959
$ws.Range("E4").Value = 1.033418652263601
$ws.Range("I4").Value = 1.038094684482384
$ws.Range("J4").Value = 1.038778771713204
$ws.Range("K4").Value = 1.040801373692023
$ws.Range("L4").Value = 1.035909162500885
$ws.Range("N4").Value = 1.016899124970719

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034495389461453
$ws.Range("D5").Value = 1.03852000027562
$ws.Range("E5").Value = 1.033636148109005
$ws.Range("I5").Value = 1.03816501718102
$ws.Range("J5").Value = 1.038937028540304
$ws.Range("K5").Value = 1.040945970620281
$ws.Range("L5").Value = 1.03607419755073
$ws.Range("N5").Value = 1.016951725309861

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034538361235888
$ws.Range("D6").Value = 1.038553031061469
$ws.Range("E6").Value = 1.033672672780054
$ws.Range("I6").Value = 1.038176807080767
$ws.Range("J6").Value = 1.038963592851842
$ws.Range("K6").Value = 1.040970238653516
$ws.Range("L6").Value = 1.036101904689855
$ws.Range("N6").Value = 1.016960554099605

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03424288432365
$ws.Range("D7").Value = 1.0383258989314
$ws.Range("E7").Value = 1.03342155804002
$ws.Range("I7").Value = 1.038095625562163
$ws.Range("J7").Value = 1.038780886864436
$ws.Range("K7").Value = 1.040803306501327
$ws.Range("L7").Value = 1.035911367905794
$ws.Range("N7").Value = 1.01689982802241

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033007972275285
$ws.Range("D8").Value = 1.037376377558086
$ws.Range("E8").Value = 1.032372839606065
$ws.Range("I8").Value = 1.037753515144601
$ws.Range("J8").Value = 1.03801609638793
$ws.Range("K8").Value = 1.040104054542059
$ws.Range("L8").Value = 1.03511453157923
$ws.Range("N8").Value = 1.016645565414665

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.030832454660207
$ws.Range("D9").Value = 1.03570273003741
$ws.Range("E9").Value = 1.030528337218146
$ws.Range("I9").Value = 1.037140212144237
$ws.Range("J9").Value = 1.03666431317452
$ws.Range("K9").Value = 1.038866277852982
$ws.Range("L9").Value = 1.033708868188932
$ws.Range("N9").Value = 1.016195894328354

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.0293825287939
$ws.Range("D10").Value = 1.034586721567955
$ws.Range("E10").Value = 1.029301075559519
$ws.Range("I10").Value = 1.036724343015049
$ws.Range("J10").Value = 1.035760401194268
$ws.Range("K10").Value = 1.038037377230719
$ws.Range("L10").Value = 1.032770780592467
$ws.Range("N10").Value = 1.01589503799665

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028754803373323
$ws.Range("D11").Value = 1.034103434714151
$ws.Range("E11").Value = 1.028770243202326
$ws.Range("I11").Value = 1.036542613317072
$ws.Range("J11").Value = 1.035368361596305
$ws.Range("K11").Value = 1.037677583421868
$ws.Range("L11").Value = 1.032364357235922
$ws.Range("N11").Value = 1.015764512626725

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02852165418596
$ws.Range("D12").Value = 1.033923914430129
$ws.Range("E12").Value = 1.028573156720011
$ws.Range("I12").Value = 1.036474862565895
$ws.Range("J12").Value = 1.035222645287022
$ws.Range("K12").Value = 1.037543809432316
$ws.Range("L12").Value = 1.032213360514051
$ws.Range("N12").Value = 1.015715992037946

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028571664744337
$ws.Range("D13").Value = 1.033962422394702
$ws.Range("E13").Value = 1.028615428425791
$ws.Range("I13").Value = 1.036489406557995
$ws.Range("J13").Value = 1.035253906232354
$ws.Range("K13").Value = 1.037572510301737
$ws.Range("L13").Value = 1.032245751302788
$ws.Range("N13").Value = 1.015726401568181

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028735530858934
$ws.Range("D14").Value = 1.034088595627769
$ws.Range("E14").Value = 1.0287539501645
$ws.Range("I14").Value = 1.036537018080241
$ws.Range("J14").Value = 1.03535631858134
$ws.Range("K14").Value = 1.037666528280617
$ws.Range("L14").Value = 1.032351876467901
$ws.Range("N14").Value = 1.015760502666164

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028836496317872
$ws.Range("D15").Value = 1.034166334356652
$ws.Range("E15").Value = 1.028839309713094
$ws.Range("I15").Value = 1.036566320229577
$ws.Range("J15").Value = 1.035419405598117
$ws.Range("K15").Value = 1.037724438585714
$ws.Range("L15").Value = 1.032417259338381
$ws.Range("N15").Value = 1.015781508497291

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029424191068334
$ws.Range("D16").Value = 1.034618794816019
$ws.Range("E16").Value = 1.029336317473464
$ws.Range("I16").Value = 1.036736368976261
$ws.Range("J16").Value = 1.03578640615807
$ws.Range("K16").Value = 1.038061237208507
$ws.Range("L16").Value = 1.032797748899062
$ws.Range("N16").Value = 1.015903695238884

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029792864239515
$ws.Range("D17").Value = 1.034902599322606
$ws.Range("E17").Value = 1.029648233413777
$ws.Range("I17").Value = 1.036842593150647
$ws.Range("J17").Value = 1.036016445177521
$ws.Range("K17").Value = 1.03827226856891
$ws.Range("L17").Value = 1.033036360245023
$ws.Range("N17").Value = 1.015980272335551

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030007914865167
$ws.Range("D18").Value = 1.035068132998219
$ws.Range("E18").Value = 1.02983022459328
$ws.Range("I18").Value = 1.036904392138256
$ws.Range("J18").Value = 1.036150561268052
$ws.Range("K18").Value = 1.03839527517254
$ws.Range("L18").Value = 1.03317551641728
$ws.Range("N18").Value = 1.016024914069186

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03008124321394
$ws.Range("D19").Value = 1.035124574872833
$ws.Range("E19").Value = 1.029892288327205
$ws.Range("I19").Value = 1.036925436880334
$ws.Range("J19").Value = 1.036196280887136
$ws.Range("K19").Value = 1.038437202901933
$ws.Range("L19").Value = 1.033222961354307
$ws.Range("N19").Value = 1.016040131600558

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029753308078651
$ws.Range("D20").Value = 1.034872150261815
$ws.Range("E20").Value = 1.029614761996174
$ws.Range("I20").Value = 1.03683121282487
$ws.Range("J20").Value = 1.035991770552155
$ws.Range("K20").Value = 1.038249635626995
$ws.Range("L20").Value = 1.033010761753355
$ws.Range("N20").Value = 1.015972058858559

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028687275925981
$ws.Range("D21").Value = 1.034051440912551
$ws.Range("E21").Value = 1.028713156518574
$ws.Range("I21").Value = 1.036523004522613
$ws.Range("J21").Value = 1.035326163323011
$ws.Range("K21").Value = 1.037638845926431
$ws.Range("L21").Value = 1.032320626155793
$ws.Range("N21").Value = 1.015750461779849

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028017111999195
$ws.Range("D22").Value = 1.033535394063822
$ws.Range("E22").Value = 1.028146792757104
$ws.Range("I22").Value = 1.036327785675091
$ws.Range("J22").Value = 1.03490711837075
$ws.Range("K22").Value = 1.037254063284795
$ws.Range("L22").Value = 1.031886520055708
$ws.Range("N22").Value = 1.015610917187935

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028372369562502
$ws.Range("D23").Value = 1.033808963128862
$ws.Range("E23").Value = 1.028446984047422
$ws.Range("I23").Value = 1.036431410838167
$ws.Range("J23").Value = 1.035129314026548
$ws.Range("K23").Value = 1.037458115128962
$ws.Range("L23").Value = 1.032116665683438
$ws.Range("N23").Value = 1.01568491295491

$ws.Range("B24").Value = 1.019999999999999
$ws.Range("C24").Value = 1.029771181774773
$ws.Range("D24").Value = 1.034885908896312
$ws.Range("E24").Value = 1.029629886116436
$ws.Range("I24").Value = 1.036836355597929
$ws.Range("J24").Value = 1.036002920144173
$ws.Range("K24").Value = 1.038259862740111
$ws.Range("L24").Value = 1.033022328677246
$ws.Range("N24").Value = 1.015975770250651

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031394807586021
$ws.Range("D25").Value = 1.036135456224255
$ws.Range("E25").Value = 1.031004766322295
$ws.Range("I25").Value = 1.037300001945384
$ws.Range("J25").Value = 1.037014265171065
$ws.Range("K25").Value = 1.039186931637363
$ws.Range("L25").Value = 1.034072442314398
$ws.Range("N25").Value = 1.016312336046204
